# Rename "jyothi" -> "jyothiaaS" in B2 and "jyothia" -> "jyothiaaS" in B3
# (both rows converge on the same UserName value), then move the
# active selection to B2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "jyothiaaS"
$ws.Range("B3").Value = "jyothiaaS"

[void]$ws.Range("B2").Select()
